# Apply weekly update: insert two new price records (rows 226-227) for
# "Betarraga" at "Vega Modelo de Temuco", pushing the existing historical
# rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 226 (existing rows 226+ shift down by 2)
$ws.Range("A226:A227").EntireRow.Insert()

# New row 226
$ws.Cells.Item(226, 1).Value2  = 10
$ws.Cells.Item(226, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(226, 3).Value2  = "La Araucanía"
$ws.Cells.Item(226, 4).Value2  = 44508
$ws.Cells.Item(226, 5).Value2  = 9
$ws.Cells.Item(226, 6).Value2  = 100114014
$ws.Cells.Item(226, 7).Value2  = "Betarraga"
$ws.Cells.Item(226, 8).Value2  = "Sin especificar"
$ws.Cells.Item(226, 9).Value2  = "Primera"
$ws.Cells.Item(226, 10).Value2 = 70
$ws.Cells.Item(226, 11).Value2 = 8000
$ws.Cells.Item(226, 12).Value2 = 9000
$ws.Cells.Item(226, 13).Value2 = 8429
$ws.Cells.Item(226, 14).Value2 = "$/docena de paquetes"
$ws.Cells.Item(226, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(226, 16).Value2 = 702
$ws.Cells.Item(226, 17).Value2 = 12
$ws.Cells.Item(226, 18).Value2 = "Hortaliza"

# New row 227
$ws.Cells.Item(227, 1).Value2  = 10
$ws.Cells.Item(227, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(227, 3).Value2  = "La Araucanía"
$ws.Cells.Item(227, 4).Value2  = 44508
$ws.Cells.Item(227, 5).Value2  = 9
$ws.Cells.Item(227, 6).Value2  = 100114014
$ws.Cells.Item(227, 7).Value2  = "Betarraga"
$ws.Cells.Item(227, 8).Value2  = "Sin especificar"
$ws.Cells.Item(227, 9).Value2  = "Primera"
$ws.Cells.Item(227, 10).Value2 = 100
$ws.Cells.Item(227, 11).Value2 = 8000
$ws.Cells.Item(227, 12).Value2 = 8000
$ws.Cells.Item(227, 13).Value2 = 8000
$ws.Cells.Item(227, 14).Value2 = "$/docena de paquetes"
$ws.Cells.Item(227, 15).Value2 = "Región del Maule"
$ws.Cells.Item(227, 16).Value2 = 667
$ws.Cells.Item(227, 17).Value2 = 12
$ws.Cells.Item(227, 18).Value2 = "Hortaliza"
